# Apply updated crypto price/volume data (and 3 coin-row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.553.61"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.312.32"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'520.14"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").Value = "'131.74"
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.533"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("D9").Value = "2.330.05"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'23.44"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.728.81"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "56.556.58"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "2.331.49"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "'334.36"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("D20").Value = "'10.37"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").Value = "'4.14"
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("D22").Value = "'6.76"
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "'61.38"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").Value = "'8.70"
$ws.Range("E25").Value = "  +7.99%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("E28").Value = "  +2.52%  "
$ws.Range("D29").Value = "'170.50"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").Value = "0.0₃0718"
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("D32").Value = "'6.08"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").Value = "'18.35"
$ws.Range("E33").Value = "  -0.76%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'0.994"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").Value = "'1.25"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").Value = "'0.895"
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'3.93"
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").Value = "'38.80"
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("D41").Value = "'148.31"
$ws.Range("E41").Value = "  +5.67%  "
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("D43").Value = "'286.76"
$ws.Range("E43").Value = "  +3.31%  "
$ws.Range("D44").Value = "'3.58"
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").Value = "'5.08"
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("D46").Value = "'0.0926"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").Value = "'0.0499"
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("D48").Value = "'0.557"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "'18.31"
$ws.Range("E49").Value = "  +2.69%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0214"
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("B51").Value = "Polygon"
$ws.Range("C51").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D51").Value = "'0.377"
$ws.Range("E51").Value = "  -0.57%  "
